# PM15 Tidsregistrering for Laila.xlsx - add four new time-registration rows
# (rows 15-18) with role/date/start/end entries, mirroring the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 15: UI Design til UC08 / User-Interface Designer ---
$ws.Cells.Item(15, 1).Value = "UI Design til UC08"
$ws.Cells.Item(15, 2).Value = "User-Interface Designer"
$ws.Cells.Item(15, 3).Value = 43894
$ws.Cells.Item(15, 4).Value = 0.35416666666666669
$ws.Cells.Item(15, 5).Value = 0.44097222222222227

# --- Row 16: Grund UI til UC08 / User-Interface Designer ---
$ws.Cells.Item(16, 1).Value = "Grund UI til UC08"
$ws.Cells.Item(16, 2).Value = "User-Interface Designer"
$ws.Cells.Item(16, 3).Value = 43894
$ws.Cells.Item(16, 4).Value = 0.44097222222222227
$ws.Cells.Item(16, 5).Value = 0.5625

# --- Row 18 entered before row 17 (matches the author's original shared-string
# insertion order: "AD09" lands before "Kundemøde" in sharedStrings.xml) ---
# Row 18: AD09 / Requirement Specifier
$ws.Cells.Item(18, 1).Value = "AD09"
$ws.Cells.Item(18, 2).Value = "Requirement Specifier"
$ws.Cells.Item(18, 3).Value = 43894
$ws.Cells.Item(18, 4).Value = 0.61805555555555558
$ws.Cells.Item(18, 5).Value = 0.66666666666666663

# Row 17: Kundemøde / Requirement Specifier
$ws.Cells.Item(17, 1).Value = "Kundemøde"
$ws.Cells.Item(17, 2).Value = "Requirement Specifier"
$ws.Cells.Item(17, 3).Value = 43894
$ws.Cells.Item(17, 4).Value = 0.59375
$ws.Cells.Item(17, 5).Value = 0.61458333333333337

# Move the active selection to E19, matching the saved view state.
$ws.Range("E19").Select()
